# pubfinder_list_NLA.xlsx — "Zu den Linden Utzmannsbach added"
#
# Inserts one new data row (new row 100) into the query-table range of the
# "pubfinder_list_NLA_ok" sheet, pushing the previous rows 100-104 down to
# 101-105, then fills the new row with the "Zu den Linden" pub entry,
# wires up its Mail/Homepage hyperlinks, and extends the table / defined
# name / dimension to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pubfinder_list_NLA_ok")

# --- 1. Insert a new row at 100, shifting existing rows 100:104 -> 101:105
$ws.Range("A100").EntireRow.Insert()

# --- 2. Fill the new row 100 with the new pub's data
$ws.Range("A100").Value = "Zu den Linden"
$ws.Range("B100").Value = "fränkisch"
$ws.Range("C100").Value = "Utzmannsbach 11"
$ws.Range("D100").Value = 91245
$ws.Range("E100").Value = "Simmelsdorf"
$ws.Range("F100").Value = "+499155 446"
$ws.Range("H100").Value = "https://www.gasthaus-zu-den-linden.de/"
$ws.Range("G100").Value = "info@gasthaus-zu-den-linden.de"
$ws.Range("I100").Value = "x"
$ws.Range("J100").Value = "x"
$ws.Range("K100").Value = "Mi"
$ws.Range("L100").Value = "Do"
$ws.Range("M100").Value = "Fr"
$ws.Range("N100").Value = "Sa"
$ws.Range("O100").Value = "So"
$ws.Range("P100").Value = "x"
$ws.Range("Q100").Value = "x"
$ws.Range("R100").Value = "10:00 - 15:00"
$ws.Range("S100").Value = "11:00 - 15:00"
$ws.Range("T100").Value = "11:00 - 15:00"
$ws.Range("U100").Value = "11:00 - 15:00"
$ws.Range("V100").Value = "11:00 - 15:00"
$ws.Range("W100").Value = "27"
$ws.Range("W100").Value = ""

# --- 3. Phone number keeps its text/quote-prefix look (matches other
#        "+49..." phone cells in the sheet, e.g. F67)
$ws.Range("F67").Copy()
$ws.Range("F100").PasteSpecial(-4122)

# --- 4. Wire up the Mail / Homepage hyperlinks for the new row
$ws.Hyperlinks.Add($ws.Range("H100"), "https://www.gasthaus-zu-den-linden.de/")
$ws.Hyperlinks.Add($ws.Range("G100"), "mailto:info@gasthaus-zu-den-linden.de")

# Re-apply the same look the other Mail/Homepage hyperlink cells use
# (e.g. G67/H67) instead of Excel's freshly-minted hyperlink style.
$ws.Range("G67").Copy()
$ws.Range("G100").PasteSpecial(-4122)
$ws.Range("H67").Copy()
$ws.Range("H100").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. Grow the query table / autofilter to the new extent
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:W105"))

# --- 6. Point the hidden "ExterneDaten_1" query defined name at the new range
$nm = $wb.Names.Item("pubfinder_list_NLA_ok!ExterneDaten_1")
$nm.RefersTo = "=pubfinder_list_NLA_ok!`$A`$1:`$W`$105"

# --- 7. Leave the selection where the author left off
$ws.Activate()
$ws.Range("V100").Select()
